$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A12: apply a localized (XOF) currency number format - decimal point & thousands separator
$ws.Range("A12").Value = 2341234
$ws.Range("A12").NumberFormat = '[$XOF]\ #,##0.00_);\([$XOF]\ #,##0.00\)'

# A13: apply a localized (F CFA) currency number format - decimal point & thousands separator
$ws.Range("A13").Value = 2341234
$ws.Range("A13").NumberFormat = '#,##0.00\ [$F CFA-340C];\-#,##0.00\ [$F CFA-340C]'

# Leave the active selection on the last-edited cell
$ws.Range("A13").Select() | Out-Null
